$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow

# Insert a new boolean-style row above row 8 (pushing Upload/header/Talent rows down by one),
# matching the formatting of row 7 ("Ref"), and label it "Force".
$ws.Rows("8:8").Insert()
$ws.Range("A7:G7").Copy($ws.Range("A8:G8"))
$ws.Range("A8").Value = "Force"

# Re-establish the frozen pane one row lower (it now needs to freeze through the new
# row 10 header), and restore the previous pane/selection semantics.
$win.FreezePanes = $false
$ws.Range("A11").Select()
$win.FreezePanes = $true
$ws.Range("A9").Select()
